$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 2.7
$ws.Range("J2").Value = 2.98
$ws.Range("G3").Value = 3.15
$ws.Range("I3").Value = 3.1
$ws.Range("Q3").Value = 1.86
$ws.Range("G4").Value = 3.05
$ws.Range("H4").Value = 2.62
$ws.Range("I4").Value = 3.25
$ws.Range("J4").Value = 3.45
$ws.Range("K4").Value = 4.6
$ws.Range("J5").Value = 5.3
$ws.Range("P5").Value = 2.36
$ws.Range("Q5").Value = 1.69
$ws.Range("R5").Value = 1.53
$ws.Range("S5").Value = 2.72
$ws.Range("T5").Value = 1.96
$ws.Range("X5").Value = 23
$ws.Range("AB5").Value = 30
$ws.Range("AE5").Value = 17
$ws.Range("AO5").Value = 5.9
$ws.Range("H6").Value = 1.89
$ws.Range("I6").Value = 1.91
$ws.Range("R6").Value = 1.54
$ws.Range("T6").Value = 1.69
$ws.Range("AJ6").Value = 100
$ws.Range("AK6").Value = 1000
$ws.Range("AL6").Value = 60
$ws.Range("AN6").Value = 42
$ws.Range("G8").Value = 13.5
$ws.Range("AG8").Value = 55
$ws.Range("H9").Value = 2.94
$ws.Range("N9").Value = 3.9
$ws.Range("S9").Value = 3.55
$ws.Range("T9").Value = 1.75
$ws.Range("U9").Value = 2.22
$ws.Range("Y9").Value = 12
$ws.Range("Z9").Value = 20
$ws.Range("AJ9").Value = 42
$ws.Range("AK9").Value = 28
$ws.Range("AL9").Value = 48
$ws.Range("AO9").Value = 36
$ws.Range("S10").Value = 2.7
$ws.Range("T10").Value = 1.64
$ws.Range("U10").Value = 2.44
$ws.Range("Y10").Value = 12
$ws.Range("AA10").Value = 23
$ws.Range("AE10").Value = 19
$ws.Range("AJ10").Value = 80
$ws.Range("AM10").Value = 70
$ws.Range("AN10").Value = 36
$ws.Range("AO10").Value = 10
$ws.Range("I11").Value = 1.74
$ws.Range("N11").Value = 6
$ws.Range("P11").Value = 2.7
$ws.Range("Q11").Value = 1.55
$ws.Range("R11").Value = 1.7
$ws.Range("S11").Value = 2.34
$ws.Range("U11").Value = 2.52
$ws.Range("X11").Value = 27
$ws.Range("Y11").Value = 13.5
$ws.Range("AA11").Value = 19
$ws.Range("AB11").Value = 26
$ws.Range("AC11").Value = 11
$ws.Range("AE11").Value = 1000
$ws.Range("AF11").Value = 46
$ws.Range("AI11").Value = 25
$ws.Range("AK11").Value = 55
$ws.Range("AN11").Value = 44
$ws.Range("AO11").Value = 6.8
